# ControlAndPoseEstimationAnalysis.xlsx - "Se agrego archivos de Rodrigo y Jairo"
#
# The substantive change in this revision is the step-time input cell
# (P2!H2, labelled "Tiempo escalon lanzar") being corrected from 1.9 to
# 1.85. Every other numeric difference in the sheet (K50, K51, K90, K91,
# K125, K126, K164, K165, K44, K47, K84, K87, K119, K122, K158, K161, ...)
# is a formula that depends on H2, so updating H2 and letting Excel
# recalculate reproduces all of them automatically.
#
# The revision also left the cursor/selection on H2 (where the edit was
# made) instead of its previous position (K47), so we update the active
# selection on sheet P2 to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("P2")

# Core data edit: step time changed from 1.9 to 1.85
$ws.Range("H2").Value = 1.85

# Leave the selection where the author left it after making the edit
$ws.Activate()
$ws.Range("H2").Select()
